# finalize the import export excel of user feature
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("dob") switches from free-text dates to real date values,
#     formatted with a new custom number format (yyyy/mm/dd). ---
$ws.Range("C2:C98").NumberFormat = "[$-1010000]yyyy/mm/dd;@"
$ws.Range("C2:C98").HorizontalAlignment = -4108
$ws.Range("C2:C98").VerticalAlignment = -4108

$ws.Range("C2").Value = 36598
$ws.Range("C3").Value = 34113
$ws.Range("C4").Value = 36598
$ws.Range("C5").Value = 34113
$ws.Range("C6").Value = 36598
$ws.Range("C7").Value = 36598
$ws.Range("C8").Value = 34113

# Rows 9-98: blank placeholder rows in column C, carrying the new date format
# so admins can paste/import additional rows (C9 already existed; C10:C98 are new).

# --- Selection moves to D10 (where the user's cursor ended up) ---
$ws.Range("D10").Select()
